$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing parameter bounds (bmin / bmax columns) ---

# PAct1_LacI (row 2): bmax 3 -> 100
$ws.Range("C2").Value = 100

# P4Lacn_cit (row 3): bmin 15 -> 0.1, bmax 25 -> 100
$ws.Range("B3").Value = 0.1
$ws.Range("C3").Value = 100

# dLacI (row 4): bmin 1E-4 -> 1E-5
$ws.Range("B4").Value = 0.00001

# dCit (row 5): bmin 1E-4 -> 1E-5
$ws.Range("B5").Value = 0.00001

# LacI_rep_WT (row 6): bmax 0.05 -> 10
$ws.Range("C6").Value = 10

# nLacI (row 8): bmin 1 -> 0.5, bmax 5 -> 20
$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = 20

# LacI_rep_W220F (row 12): bmax 0.01 -> 1
$ws.Range("C12").Value = 1

# P_4Lacn_LacI (row 15): bmin 1 -> 0.1, bmax 10 -> 100
$ws.Range("B15").Value = 0.1
$ws.Range("C15").Value = 100

# LacI_rep_3mut (row 17): bmax 1E-3 -> 0.1
$ws.Range("C17").Value = 0.1

# pt7_LacI (row 19): bmin 1 -> 0.1, bmax 10 -> 100
$ws.Range("B19").Value = 0.1
$ws.Range("C19").Value = 100

# P3_Lacn_5_cit (row 20): bmin 5 -> 0.1, bmax 25 -> 100
$ws.Range("B20").Value = 0.1
$ws.Range("C20").Value = 100

# P3_Lacn_5_cit_L (row 21): bmin 1E-3 -> 1E-6
$ws.Range("B21").Value = 0.000001

# dLacI_pt7 (row 22): bmin 1 -> 1E-3
$ws.Range("B22").Value = 0.001

# nLacI_P3 (row 23): bmin 1 -> 0.5, bmax 3 -> 20
$ws.Range("B23").Value = 0.5
$ws.Range("C23").Value = 20

# --- Add new parameter row 24: LacI_rep_3mut_P3 ---
# Match the surrounding rows' font (Calibri 11) for columns A-F so the new
# row doesn't fall back to the workbook default (12pt) look.
$ws.Range("A24:F24").Font.Size = 11

$ws.Range("G24").Value = "theta_{LacI_W220F_Q60G_T167A-pt7}"
$ws.Range("A24").Value = "LacI_rep_3mut_P3"
$ws.Range("B24").Value = 0.00001
$ws.Range("C24").Value = 0.1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = "yes"

# --- Update view state (scroll position / selection) ---
# Excel scrolled the sheet down (topLeftCell -> A3) and the active selection
# moved to the new row just past the appended data (A25).
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A25").Select() | Out-Null
